# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Ciruela" (Terminal Hortofrutícola
# Agro Chillán) above the current row 81, pushing the existing rows 81-88
# down to 83-90 and extending the used range to A1:T90.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 81:88 down by two rows, duplicating formatting (style) from
# the row above, just like Excel's native "Insert Rows" command.
$ws.Rows("81:82").Insert()

# New row 81: Ciruela Angeleno, Primera
$ws.Range("A81").Value = 7
$ws.Range("B81").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C81").Value = "Ñuble"
$ws.Range("D81").Value = 44995
$ws.Range("E81").Value = 16
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100103
$ws.Range("H81").Value = "Frutos de hueso (carozo)"
$ws.Range("I81").Value = 100103002
$ws.Range("J81").Value = "Ciruela"
$ws.Range("K81").Value = "Angeleno"
$ws.Range("L81").Value = "Primera"
$ws.Range("M81").Value = 80
$ws.Range("N81").Value = 10000
$ws.Range("O81").Value = 10000
$ws.Range("P81").Value = 10000
$ws.Range("Q81").Value = "$/bandeja 18 kilos granel"
$ws.Range("R81").Value = "Región de O'Higgins"
$ws.Range("S81").Value = 556
$ws.Range("T81").Value = 18

# New row 82: Ciruela Larry Ann, Primera
$ws.Range("A82").Value = 7
$ws.Range("B82").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C82").Value = "Ñuble"
$ws.Range("D82").Value = 44995
$ws.Range("E82").Value = 16
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100103
$ws.Range("H82").Value = "Frutos de hueso (carozo)"
$ws.Range("I82").Value = 100103002
$ws.Range("J82").Value = "Ciruela"
$ws.Range("K82").Value = "Larry Ann"
$ws.Range("L82").Value = "Primera"
$ws.Range("M82").Value = 50
$ws.Range("N82").Value = 9000
$ws.Range("O82").Value = 9000
$ws.Range("P82").Value = 9000
$ws.Range("Q82").Value = "$/bandeja 18 kilos granel"
$ws.Range("R82").Value = "Región de O'Higgins"
$ws.Range("S82").Value = 500
$ws.Range("T82").Value = 18
